$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: "What is the interquartile range..." question type changes
# from "mult choice" to "single choice"
$ws.Range("E4").Value = "single choice"

# Row 16: fill in the previously-empty Question / Name of file / Type of Sol
# columns for the "assumptions of matched-pairs t-test" learning objective
# (breaking up old 2 sample ttest questions)
$ws.Range("C16").Value = "read test-statistic, read sided t-test"
$ws.Range("D16").Value = "schoice-2samtt-interpret-t, schoice-2samtt-interpret-alt"
$ws.Range("E16").Value = "single choice"

# Update the active cell selection to match the authored state
$ws.Range("D17").Select()
